# Add the new row of profit data produced by the 2025-09-06 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (matching the existing rows,
# which store "MM/DD/YYYY" strings rather than real date values). Force
# the cell to Text first so Excel doesn't auto-convert the string into a
# date serial, then clear the formatting again so the cell is left with
# the default (unstyled) look, just like its neighbours.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "09/06/2025"
$ws.Range("A20").ClearFormats()

# Column B holds the numeric profit value for the day.
$ws.Range("B20").Value = 12477.94
